# fix challenges table problem
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the now-unused trailing empty rows (8-14), which only carried
# formatting (style index 2 on column A).
$ws.Range("A8:A14").EntireRow.Delete() | Out-Null

# Add the new "status" column header + values.
$ws.Range("E1").Value = "status"
$ws.Range("E2").Value = 0
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 10
$ws.Range("E7").Value = 1000

# Column width adjustments (closest settable values given the host's
# internal 1/6-character rounding grid for ColumnWidth).
$ws.Columns.Item(4).ColumnWidth = 10
$ws.Columns.Item(5).ColumnWidth = 15.833333333333334
$ws.Columns.Item(6).ColumnWidth = 19.833333333333332

# Selection moves to C18 in the final saved state.
$ws.Range("C18").Select()
